# The post "「間違いで賢くなる、痛みで強くなる」" (row 478) was removed from the
# blog/posts sheet. Deleting the entire row shifts every subsequent row
# up by one (and Excel keeps the sheet dimension in sync automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(478).Delete()
